$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the J/K (SPID/CS0) columns as "does not work" by striking through
# the relevant header/legend cells (same visual treatment already used for
# column J7/K7 "D18"/"D23").
$ws.Range("J2:K3").Font.Strikethrough = $true
$ws.Range("J4:K4").Font.Strikethrough = $true
$ws.Range("J5:K5").Font.Strikethrough = $true
$ws.Range("J6:K6").Font.Strikethrough = $true

# --- Update/expand the IC2 and SPI interface descriptions (still on rows
# 18/19 at this point, before the new row is inserted below).
$ws.Range("B18").Value = "The IC2 interface used by some components. Example: BMP180 Atmospheric Pressure Sensor"
$ws.Range("B19").Value = "The SPI interface used by some components. Example: DS3231 RTC Module"

# --- Insert a new "TOUCH" row into the Features legend (between "ADC" row
# and the IC2/SPI/Espruino rows), pushing the following rows down by one.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = "TOUCH"
$ws.Range("B18").Value = "Not sure if it works in Espruino: http://forum.espruino.com/conversations/328396/"

# --- Move the hyperlink that used to live on B20 down to B21 (its new row).
$ws.Range("B21").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B21"), "https://github.com/espruino/Espruino/issues/1574") | Out-Null
$ws.Range("B21").HorizontalAlignment = -4131

# --- Leave the selection where the author last left it before saving.
$ws.Range("B22").Select() | Out-Null
